$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header labels (row 1) to reflect new flow node numbering ---
$ws.Range("B1").Value = "F_2_3_t"
$ws.Range("C1").Value = "F_3_4_t"
$ws.Range("D1").Value = "F_4_0_t"
$ws.Range("E1").Value = "F_4_5_t"
$ws.Range("F1").Value = "F_4_7_t"
$ws.Range("G1").Value = "F_5_6_t"
$ws.Range("H1").Value = "F_5_7_t"
$ws.Range("I1").Value = "F_6_0_t"
$ws.Range("J1").Value = "F_6_1_t"
$ws.Range("K1").Value = "F_7_0_t"
$ws.Range("L1").Value = "F_7_1_t"
$ws.Range("M1").Value = "F_1_2_t"

# --- Update recalculated simulation values (columns L, M, N) ---
$ws.Range("M52").Value = 0.5102828806193229
$ws.Range("M53").Value = 0.516669718407005
$ws.Range("M54").Value = 0.5250598450728583
$ws.Range("M55").Value = 0.5549266916070887
$ws.Range("M56").Value = 0.5606662037032415
$ws.Range("M57").Value = 0.6232718818418949
$ws.Range("M58").Value = 0.6335318285140253
$ws.Range("M59").Value = 0.5955934837702366
$ws.Range("M60").Value = 0.5983913130152423
$ws.Range("M61").Value = 0.6137176582650126
$ws.Range("M62").Value = 0.6277536792345876
$ws.Range("M63").Value = 0.6359854366096307
$ws.Range("M64").Value = 0.6503700033180759
$ws.Range("M65").Value = 0.7280687666642501
$ws.Range("M66").Value = 0.7570927958566408
$ws.Range("M67").Value = 0.8166769708856849
$ws.Range("M68").Value = 0.8255399075291283
$ws.Range("M69").Value = 0.8221884962613993
$ws.Range("M70").Value = 0.8896522866220384
$ws.Range("M71").Value = 0.915464558246726
$ws.Range("M72").Value = 0.9035809647651953
$ws.Range("M73").Value = 0.9959184452730141
$ws.Range("M74").Value = 1.045956681996581
$ws.Range("M75").Value = 1.092495866739968
$ws.Range("M76").Value = 0.9803243046049885
$ws.Range("M77").Value = 0.9845418964038379
$ws.Range("M78").Value = 1.10172151779957
$ws.Range("M79").Value = 1.219511167928284
$ws.Range("M80").Value = 1.31484762270926
$ws.Range("M81").Value = 1.36311253422606
$ws.Range("M82").Value = 1.279229416176819
$ws.Range("M83").Value = 1.280037547355418
$ws.Range("L84").Value = 0.504028518323153
$ws.Range("M84").Value = 1.279940383139227
$ws.Range("N84").Value = 0.008741271316853311
$ws.Range("M85").Value = 1.419851644933041
$ws.Range("M86").Value = 1.546631001058964
$ws.Range("M87").Value = 1.712910885515025
$ws.Range("M88").Value = 1.744789858238614
$ws.Range("M89").Value = 1.876003014706268
$ws.Range("M90").Value = 2.010732449452155
$ws.Range("M91").Value = 2.107265313709849
$ws.Range("M92").Value = 2.352869598949975
$ws.Range("M93").Value = 2.418310037210902
$ws.Range("M94").Value = 2.504535692240622
$ws.Range("M95").Value = 2.553707092759789
$ws.Range("M96").Value = 2.725382211854786
$ws.Range("M97").Value = 2.871723071146878
$ws.Range("M98").Value = 3.10536315657654
$ws.Range("M99").Value = 3.326008930096381
$ws.Range("M100").Value = 3.474731694744438
$ws.Range("M101").Value = 3.79204514918327
$ws.Range("M102").Value = 4.054575995693402
$ws.Range("M103").Value = 4.010294357635344
$ws.Range("M104").Value = 4.264767908251616
$ws.Range("M105").Value = 4.514829405688729
$ws.Range("M106").Value = 4.99509157636402
$ws.Range("M107").Value = 5.508717731141357
$ws.Range("M108").Value = 5.964101143454548
$ws.Range("M109").Value = 6.398077947644984
$ws.Range("M110").Value = 6.414563993245897
$ws.Range("M111").Value = 5.801014989085774
$ws.Range("M112").Value = 7.062141365224273
$ws.Range("M113").Value = 7.445890212056876
$ws.Range("M114").Value = 8.063052196594999
$ws.Range("M115").Value = 8.578405521463791
$ws.Range("M116").Value = 9.124963293040034
$ws.Range("M117").Value = 9.502894731612617
$ws.Range("M118").Value = 10.15633386560947
$ws.Range("M119").Value = 11.01819513312171
$ws.Range("M120").Value = 11.15030017533369
$ws.Range("M121").Value = 13.51853464436738
$ws.Range("M122").Value = 13.22658022345069
$ws.Range("L123").Value = 3.643527638993817
$ws.Range("M123").Value = 13.70503296232731
$ws.Range("L124").Value = 3.766236285362266
$ws.Range("M124").Value = 15.61485869122209
$ws.Range("L125").Value = 3.888075800604792
$ws.Range("M125").Value = 16.36178729181363
$ws.Range("L126").Value = 4.010985667144613
$ws.Range("M126").Value = 17.25287387176779
$ws.Range("L127").Value = 4.137172850243049
$ws.Range("M127").Value = 18.10135676296288
$ws.Range("L128").Value = 4.268342152917543
$ws.Range("M128").Value = 19.0761957257819
$ws.Range("L129").Value = 4.405018343013344
$ws.Range("M129").Value = 20.10300639473459
$ws.Range("L130").Value = 4.546221523176349
$ws.Range("M130").Value = 21.08422671372097
$ws.Range("L131").Value = 4.689639300927146
$ws.Range("M131").Value = 22.09667840307107
$ws.Range("L132").Value = 4.83226242520357
$ws.Range("M132").Value = 23.09340491993343
$ws.Range("N132").Value = 0
$ws.Range("L133").Value = 4.971234061854747
$ws.Range("M133").Value = 22.55252506063162
$ws.Range("N133").Value = 0
$ws.Range("L134").Value = 5.104777522141549
$ws.Range("M134").Value = 23.10328568696132
$ws.Range("N134").Value = 0
$ws.Range("L135").Value = 5.232397259208467
$ws.Range("M135").Value = 23.63030746751238
$ws.Range("N135").Value = 0
$ws.Range("L136").Value = 5.354958701814573
$ws.Range("M136").Value = 24.1509965347631
$ws.Range("N136").Value = 0.02899257651740307
$ws.Range("L137").Value = 5.474565476900663
$ws.Range("M137").Value = 24.68793678712501
$ws.Range("N137").Value = 0.2186348278125633
$ws.Range("L138").Value = 5.594801192938578
$ws.Range("M138").Value = 25.2679229503968
$ws.Range("N138").Value = 0.4144811464015073
$ws.Range("L139").Value = 5.721230630560602
$ws.Range("M139").Value = 25.92341006260021
$ws.Range("N139").Value = 0.6178008467246232
$ws.Range("L140").Value = 5.861156889361621
$ws.Range("M140").Value = 26.69348991873918
$ws.Range("N140").Value = 0.8318260143913753
$ws.Range("L141").Value = 6.02109043905461
$ws.Range("M141").Value = 27.62074690740913
$ws.Range("N141").Value = 1.05933408629137
$ws.Range("L142").Value = 6.201265911305343
$ws.Range("M142").Value = 28.73513745805013
$ws.Range("N142").Value = 1.298935738521695
$ws.Range("L143").Value = 6.389015004885707
$ws.Range("M143").Value = 30.03116016780058
$ws.Range("N143").Value = 1.547177203883071
$ws.Range("L144").Value = 6.555731277362504
$ws.Range("M144").Value = 31.4390062916234
$ws.Range("N144").Value = 1.777675720261958
$ws.Range("L145").Value = 6.661597245219092
$ws.Range("M145").Value = 32.83101012412658
$ws.Range("N145").Value = 1.963439913860894
$ws.Range("L146").Value = 6.670145285603253
$ws.Range("M146").Value = 34.0587482367197
$ws.Range("N146").Value = 2.080505583751598
$ws.Range("L147").Value = 6.567681657590218
$ws.Range("M147").Value = 35.03242968957792
$ws.Range("N147").Value = 2.119902808124075
$ws.Range("L148").Value = 6.37748740627446
$ws.Range("M148").Value = 35.81029870481003
$ws.Range("N148").Value = 2.095713462210971
$ws.Range("L149").Value = 6.156114452447702
$ws.Range("M149").Value = 36.63741007721087
$ws.Range("N149").Value = 2.041650781456637
$ws.Range("L150").Value = 5.965134961981072
$ws.Range("M150").Value = 37.86712479985761
$ws.Range("N150").Value = 1.992612980593269
$ws.Range("L151").Value = 5.829373161168148
$ws.Range("M151").Value = 39.75683710949776
$ws.Range("N151").Value = 1.959130348477006
$ws.Range("L152").Value = 5.711303899882349
$ws.Range("M152").Value = 42.24500285206719
$ws.Range("N152").Value = 1.914543577730211
